# Apply a cyclic rotation of the species-record fields across rows 2-5.
# Row 2's values move to row 5, and rows 3,4,5 each shift up into 2,3,4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the rotating data for this record.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values of row 2 (rows 3-5 will simply move up into
# 2-4, and row 2's original data ends up in row 5).
$origRow2 = @{}
foreach ($col in $cols) {
    $origRow2[$col] = $ws.Range($col + "2").Value2
}

# Shift rows 3,4,5 up into rows 2,3,4.
for ($r = 2; $r -le 4; $r++) {
    $srcRow = $r + 1
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = $ws.Range($col + $srcRow).Value2
    }
}

# Place the original row 2 values into row 5.
foreach ($col in $cols) {
    $ws.Range($col + "5").Value = $origRow2[$col]
}
